$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calcul")

# ---------------------------------------------------------------------------
# 1) Header row (row 1): rearrange "Augmentation" / "Prix suggere programme R"
#    and introduce the brand-new "Augmentation R" column header.
#    Before: H1=NPS  I1=PSPR        J1=(empty)  K1=Augmentation  L1=TxAug
#    After : H1=NPS  I1=Augmentation J1=PSPR     K1=Augmentation R L1=TxAug
# ---------------------------------------------------------------------------
$i1Text = $ws.Range("I1").Value2
$k1Text = $ws.Range("K1").Value2

# Style: old K1 (plain bold, no fill) moves to I1; old I1 (bold + green
# fill) moves to K1 (the new "Augmentation R" column keeps the highlighted
# look); J1 picks up the bold+fill look too (it now hosts the header for
# the manual-value column). Stage old I1's format via an unused cell (X1)
# so it survives the I1 overwrite below.
$ws.Range("I1").Copy() | Out-Null
$ws.Range("X1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (stash)
$ws.Range("K1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("X1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("X1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("X1").Clear() | Out-Null
$excel.CutCopyMode = 0

$ws.Range("I1").Value = $k1Text
$ws.Range("J1").Value = $i1Text
$ws.Range("K1").Value = "Augmentation R"

# ---------------------------------------------------------------------------
# 2) Data rows 2..19: columns I/J/K/L get rearranged.
#    Before: I=manual value   J=(I-G)         K=(H-G)          L=(K/G)
#    After : I=(H-G)          J=manual value  K=(J-G)          L=(I/G)
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 19; $r++) {
    $manualVal = $ws.Range("I$r").Value2

    # Move the "manual value" look (green fill, currency) from I to J.
    $ws.Range("I$r").Copy() | Out-Null
    $ws.Range("J$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    # Move the "computed" look (plain currency) from K to I.
    $ws.Range("K$r").Copy() | Out-Null
    $ws.Range("I$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $excel.CutCopyMode = 0

    $ws.Range("J$r").Value = $manualVal
    $ws.Range("I$r").Formula = "=H$r-G$r"
    $ws.Range("K$r").Formula = "=J$r-G$r"
    $ws.Range("L$r").Formula = "=I$r/G$r"
}

# ---------------------------------------------------------------------------
# 3) Data correction on row 19: C19 2.3 -> 1.98 (ripples through the row's
#    formulas automatically on recalculation).
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 1.98

# ---------------------------------------------------------------------------
# 4) Conditional formatting range: was J2:L19, becomes K2:L19 plus I2:I19.
# ---------------------------------------------------------------------------
$oldCF = $ws.Range("J2:L19").FormatConditions
$cfLess = $oldCF.Item(1)
$cfMore = $oldCF.Item(2)
$klRange = $ws.Range("K2:L19")
$cfLess.ModifyAppliesToRange($klRange) | Out-Null
$cfMore.ModifyAppliesToRange($klRange) | Out-Null

$iRange = $ws.Range("I2:I19")
$cfILess = $iRange.FormatConditions.Add(1, 6, "0")
$cfILess.Font.Color = 393372
$cfILess.Interior.Color = 13551615
$cfIMore = $iRange.FormatConditions.Add(1, 5, "0")
$cfIMore.Font.Color = 24832
$cfIMore.Interior.Color = 13561798

# ---------------------------------------------------------------------------
# 5) Selection moves from H19 to A19.
# ---------------------------------------------------------------------------
$ws.Range("A19").Select() | Out-Null
